$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48: CENSUS2010_BLK_BG_TRCT_SHP/ ---
$ws.Range("A48").Value = "CENSUS2010_BLK_BG_TRCT_SHP/"
$ws.Range("B48").Value = "Folder with 2010 MA census block shapefiles"
$ws.Range("C48").Value = "Mass.gov"
$ws.Range("D48").Value = "https://docs.digital.mass.gov/dataset/massgis-data-datalayers-2010-us-census"
$ws.Range("E48").Value = 43723
$ws.Range("F48").Value = "Used to merge census blocks with towns (assumption that 2010 census blocks are the same as ACS blocks - GEOID10 did match up)"

# --- Row 49: townssurvey_shp/ ---
$ws.Range("A49").Value = "townssurvey_shp/"
$ws.Range("B49").Value = "Folder with MA town survey shapefiles (state of MA divided up into 351 municipalities/towns/cities)"
$ws.Range("C49").Value = "Mass.gov"
$ws.Range("D49").Value = "https://docs.digital.mass.gov/dataset/massgis-data-community-boundaries-towns-survey-points"
$ws.Range("E49").Value = 43724
$ws.Range("F49").Value = "Used to merge census blocks with towns"

# Bold font on column A cells, matching section-header look used elsewhere in the sheet
$ws.Range("A48").Font.Bold = $true
$ws.Range("A49").Font.Bold = $true

# Date formatting for the "Date accessed/downloaded" column, copied from an existing
# date cell so the workbook reuses the same style record instead of minting a new one.
$ws.Range("E3").Copy()
$ws.Range("E48").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-set numeric values (PasteSpecial(Formats) only copies formatting, not values)
$ws.Range("E48").Value = 43723
$ws.Range("E49").Value = 43724

# --- sheet view: scrolled down with D51 selected, matching the saved view state ---
$ws.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 31
